$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add the two new table columns (perrons, reizigers) at the end of Table1
$colPerrons = $lo.ListColumns.Add()
$colReizigers = $lo.ListColumns.Add()

# Name the new columns via their header cells (this updates sharedStrings + table column names)
$colPerrons.Range.Cells(1,1).Value = "perrons"
$colReizigers.Range.Cells(1,1).Value = "reizigers"

# Bulk-write the perrons/reizigers values for all 52 data rows (rows 2-53)
$arr = New-Object 'object[,]' 52,2
$arr[0,0]=2; $arr[0,1]=5548
$arr[1,0]=6; $arr[1,1]=36131
$arr[2,0]=6; $arr[2,1]=63532
$arr[3,0]=2; $arr[3,1]=6474
$arr[4,0]=5; $arr[4,1]=56752
$arr[5,0]=2; $arr[5,1]=7012
$arr[6,0]=3; $arr[6,1]=11359
$arr[7,0]=3; $arr[7,1]=31144
$arr[8,0]=2; $arr[8,1]=5786
$arr[9,0]=2; $arr[9,1]=3962
$arr[10,0]=2; $arr[10,1]=1889
$arr[11,0]=4; $arr[11,1]=3819
$arr[12,0]=6; $arr[12,1]=9791
$arr[13,0]=3; $arr[13,1]=8967
$arr[14,0]=4; $arr[14,1]=8078
$arr[15,0]=5; $arr[15,1]=3072
$arr[16,0]=2; $arr[16,1]=6738
$arr[17,0]=8; $arr[17,1]=43362
$arr[18,0]=4; $arr[18,1]=12590
$arr[19,0]=3; $arr[19,1]=15346
$arr[20,0]=2; $arr[20,1]=4840
$arr[21,0]=4; $arr[21,1]=39943
$arr[22,0]=3; $arr[22,1]=10763
$arr[23,0]=6; $arr[23,1]=43123
$arr[24,0]=16; $arr[24,1]=229788
$arr[25,0]=4; $arr[25,1]=13157
$arr[26,0]=5; $arr[26,1]=15272
$arr[27,0]=11; $arr[27,1]=171272
$arr[28,0]=3; $arr[28,1]=3140
$arr[29,0]=2; $arr[29,1]=1607
$arr[30,0]=6; $arr[30,1]=36507
$arr[31,0]=5; $arr[31,1]=25012
$arr[32,0]=3; $arr[32,1]=12896
$arr[33,0]=6; $arr[33,1]=90108
$arr[34,0]=2; $arr[34,1]=34989
$arr[35,0]=12; $arr[35,1]=80091
$arr[36,0]=7; $arr[36,1]=20360
$arr[37,0]=5; $arr[37,1]=28246
$arr[38,0]=6; $arr[38,1]=80409
$arr[39,0]=13; $arr[39,1]=104840
$arr[40,0]=4; $arr[40,1]=26979
$arr[41,0]=4; $arr[41,1]=12592
$arr[42,0]=3; $arr[42,1]=6544
$arr[43,0]=3; $arr[43,1]=20355
$arr[44,0]=5; $arr[44,1]=8334
$arr[45,0]=3; $arr[45,1]=5740
$arr[46,0]=14; $arr[46,1]=42976
$arr[47,0]=3; $arr[47,1]=7548
$arr[48,0]=3; $arr[48,1]=6386
$arr[49,0]=11; $arr[49,1]=16064
$arr[50,0]=2; $arr[50,1]=4987
$arr[51,0]=7; $arr[51,1]=7705
$ws.Range("L2:M53").Value = $arr

# Formatting: L (perrons) reuses the existing right-aligned "general number" style
# already used by columns H:K (same look as the N/Y flag columns).
$ws.Range("L2:L53").HorizontalAlignment = -4152

# M (reizigers): thousands-separated integer format, right aligned for the data,
# with the same number format (but excel-default alignment) on the header cell.
$ws.Range("M2:M53").NumberFormat = "#,##0"
$ws.Range("M2:M53").HorizontalAlignment = -4152
$ws.Range("M1").NumberFormat = "#,##0"

# Column widths (best-effort match of the post-edit layout)
$ws.Columns("A").ColumnWidth = 16.5
$ws.Columns("B").ColumnWidth = 12.5
$ws.Columns("C").ColumnWidth = 10
$ws.Columns("D").ColumnWidth = 13.83
$ws.Columns("E").ColumnWidth = 9.5
$ws.Columns("G").ColumnWidth = 10.33
$ws.Columns("H").ColumnWidth = 8
$ws.Columns("J").ColumnWidth = 5.5
$ws.Columns("K").ColumnWidth = 7.33
$ws.Columns("L").ColumnWidth = 7.33
$ws.Columns("M").ColumnWidth = 7.33

# Sheet view: zoom + selection matches the saved workbook state
$win = $ws.Application.ActiveWindow
$win.Zoom = 90
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("O20").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
